# Update predicted values in column C on the active sheet
# (retraining fix for the "PC Sun" model)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = 0.319
    5   = 0.623
    6   = 0.827
    7   = 1.345
    8   = 1.347
    11  = 0.332
    12  = 0.056
    28  = 0.07000000000000001
    29  = 0.242
    30  = 0.3
    31  = 0.334
    32  = 0.265
    33  = 0.163
    34  = 0.118
    35  = 0.052
    36  = 0.01
    51  = 0.01
    52  = 0.137
    53  = 0.431
    54  = 0.8110000000000001
    55  = 1.036
    56  = 1.101
    76  = 0.115
    77  = 0.338
    78  = 0.614
    79  = 0.758
    80  = 0.824
    81  = 0.764
    99  = 0
    100 = 0.11
    101 = 0.368
    102 = 0.732
    103 = 0.895
    105 = 0.922
    106 = 0.603
    128 = 2.105
    129 = 1.81
    130 = 0.977
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
